# This script renames a set of header labels (row 1) that are used across
# every worksheet of the workbook. The labels are PowerFactory short-circuit
# result column names; several of them had an abbreviated "_deg"/"_mw"/
# "_bus_*" suffix that needed to be corrected to the full, accurate suffix
# ("_degree" / "_mvar" / without the superfluous "_bus").
#
# Mapping of old -> new header text (applies wherever the old text occurs
# in row 1 of any worksheet):
$renameMap = @{
    "pf_ikss_from_deg"      = "pf_ikss_from_degree"
    "pf_ikss_to_deg"        = "pf_ikss_to_degree"
    "pf_va_from_deg"        = "pf_va_from_degree"
    "pf_va_to_deg"          = "pf_va_to_degree"
    "pf_q_a_from_mw"        = "pf_q_a_from_mvar"
    "pf_q_b_from_mw"        = "pf_q_b_from_mvar"
    "pf_q_c_from_mw"        = "pf_q_c_from_mvar"
    "pf_q_a_to_mw"          = "pf_q_a_to_mvar"
    "pf_q_b_to_mw"          = "pf_q_b_to_mvar"
    "pf_q_c_to_mw"          = "pf_q_c_to_mvar"
    "pf_ikss_a_from_deg"    = "pf_ikss_a_from_degree"
    "pf_ikss_b_from_deg"    = "pf_ikss_b_from_degree"
    "pf_ikss_c_from_deg"    = "pf_ikss_c_from_degree"
    "pf_ikss_a_to_deg"      = "pf_ikss_a_to_degree"
    "pf_ikss_b_to_deg"      = "pf_ikss_b_to_degree"
    "pf_ikss_c_to_deg"      = "pf_ikss_c_to_degree"
    "pf_vm_b_from_bus_pu"   = "pf_vm_b_from_pu"
    "pf_vm_c_from_bus_pu"   = "pf_vm_c_from_pu"
    "pf_vm_a_to_bus_pu"     = "pf_vm_a_to_pu"
    "pf_vm_b_to_bus_pu"     = "pf_vm_b_to_pu"
    "pf_vm_c_to_bus_pu"     = "pf_vm_c_to_pu"
    "pf_va_a_from_bus_deg"  = "pf_va_a_from_degree"
    "pf_va_b_from_bus_deg"  = "pf_va_b_from_degree"
    "pf_va_c_from_bus_deg"  = "pf_va_c_from_degree"
    "pf_va_a_to_bus_deg"    = "pf_va_a_to_degree"
    "pf_va_b_to_bus_deg"    = "pf_va_b_to_degree"
    "pf_va_c_to_bus_deg"    = "pf_va_c_to_degree"
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $headerRow = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $usedRange.Columns.Count))
    foreach ($cell in $headerRow.Cells) {
        $text = $cell.Value2
        if ($null -ne $text -and $renameMap.ContainsKey($text)) {
            $cell.Value = $renameMap[$text]
        }
    }
}
